# NerdQuiz_TODO.xlsx edit
# Commit message: "Categories default now checked"
#
# Semantic changes applied:
#  1. CategoryActivity / "Default Checked" (row 9) status Lia -> moves from
#     "In Progress" to "Done" (the actual data change driving the commit).
#  2. Cosmetic trim of trailing spaces in two shared strings:
#       - the "Übernommen von " column header (used by all 4 tables)
#       - "Methode zum prüfen ob alle Fragen gespielt wurden " (row 27)
#  3. Lower-cased "StartGameActivity" -> "startgameactivity" inside the
#     literal/error formula in B23.
#  4. Cosmetic: rename the custom cell style "Excel Built-in Input" to
#     "Excel Built-in Excel Built-in Input" (as produced by the original
#     authoring tool on re-save).
#  5. Selection cursor moved from C10 to D9 (the cell that was edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Core data edit: CategoryActivity / Default Checked / Lia -> Done
$ws.Range("D9").Value = "Done"

# 2a. Trim trailing space on the repeated column header "Übernommen von "
$ws.Range("C1").Value = "Übernommen von"
$ws.Range("C8").Value = "Übernommen von"
$ws.Range("C15").Value = "Übernommen von"
$ws.Range("C20").Value = "Übernommen von"

# 2b. Trim trailing space on the "all questions played" task text
$ws.Range("B27").Value = "Methode zum prüfen ob alle Fragen gespielt wurden"

# 3. Lower-case StartGameActivity inside the B23 note/formula
$ws.Range("B23").Formula = "=> dialog fehlt der zugriff auf 'startgameactivity' - eigener intent zur highscoreactivity?"

# 4. Rename the custom "Excel Built-in Input" cell style
try {
    $styles = $wb.Styles
    for ($i = 1; $i -le $styles.Count; $i++) {
        $s = $styles.Item($i)
        if ($s.Name -eq "Excel Built-in Input") {
            $s.Name = "Excel Built-in Excel Built-in Input"
        }
    }
} catch {
    # Older/limited hosts may not expose a writable Styles collection;
    # this cosmetic rename is best-effort and safe to skip.
}

# 5. Move the active selection to D9, matching the saved view state
[void]$ws.Range("D9").Select()
